$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- New separator row 25 ("FRIDAY") — same look as the row 9 / row 19 separators ---
$ws.Range("A9:F9").Copy() | Out-Null
$ws.Range("A25:F25").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A9:F9").Copy() | Out-Null
$ws.Range("A25:F25").PasteSpecial(-4163) | Out-Null   # xlPasteValues

# --- New log entry row 26 — repeat of the "AV Shutdown / BC / 320" task (row 6), new date+time ---
$ws.Range("A6:F6").Copy() | Out-Null
$ws.Range("A26:F26").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A6:F6").Copy() | Out-Null
$ws.Range("A26:F26").PasteSpecial(-4163) | Out-Null   # xlPasteValues

$ws.Range("B26").Value = 42587
$ws.Range("C26").Value = "1500"

$ws.Rows.Item(26).RowHeight = 45

$ws.Activate() | Out-Null
$ws.Range("F30").Select() | Out-Null
